$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: B35 becomes text "9.30-11.30" instead of time value
$ws.Range("B35").Value = "9.30-11.30"

# Row 35: C35 text extended
$ws.Range("C35").Value = "Siirtyminen törmäyksen havaitsemisesta kontaktien aiheuttamiin voimiin, Kahden laatikon törmäystarkastelu, laatikon ja tason, ja laatikon ja pisteen törmäystarkastelu"

# Row 35: G35 new hour value
$ws.Range("G35").Value = 2

# Row height grows to fit the longer wrapped text in C35
$ws.Rows(35).RowHeight = 87

# Update selection to H35
$ws.Range("H35").Select()
